$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 48679.8
$ws.Range("J117").Value = 48679.8
$ws.Range("L117").Value = 48679.8
$ws.Range("N117").Value = -57857.8
$ws.Range("H124").Value = 46476.5
$ws.Range("J124").Value = 46476.5
$ws.Range("L124").Value = 46476.5
$ws.Range("N124").Value = -56296.5
$ws.Range("H128").Value = 46088.8
$ws.Range("J128").Value = 46088.8
$ws.Range("L128").Value = 46088.8
$ws.Range("N128").Value = -56048.8
$ws.Range("H130").Value = 47384
$ws.Range("J130").Value = 47384
$ws.Range("L130").Value = 47384
$ws.Range("N130").Value = -57424
$ws.Range("H133").Value = 38827
$ws.Range("J133").Value = 38827
$ws.Range("L133").Value = 38827
$ws.Range("N133").Value = -48947

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 38311.89
$ws.Range("J80").Value = 38311.89
$ws.Range("L80").Value = 38311.89
$ws.Range("N80").Value = -40307.89
$ws.Range("H83").Value = 38311.89
$ws.Range("J83").Value = 38311.89
$ws.Range("L83").Value = 114935.67
$ws.Range("N83").Value = -124919.67
$ws.Range("H118").Value = 49409
$ws.Range("J118").Value = 49409
$ws.Range("L118").Value = 49409
$ws.Range("N118").Value = -52723
$ws.Range("H123").Value = 36790.2
$ws.Range("J123").Value = 36790.2
$ws.Range("L123").Value = 36790.2
$ws.Range("N123").Value = -46590.2
$ws.Range("H125").Value = 45398
$ws.Range("J125").Value = 45398
$ws.Range("L125").Value = 45398
$ws.Range("N125").Value = -55238
$ws.Range("H130").Value = 48421
$ws.Range("J130").Value = 48421
$ws.Range("L130").Value = 48421
$ws.Range("N130").Value = -58461
$ws.Range("H131").Value = 44136.25
$ws.Range("J131").Value = 44136.25
$ws.Range("L131").Value = 44136.25
$ws.Range("N131").Value = -54216.25
$ws.Range("H138").Value = 52800
$ws.Range("J138").Value = 52800
$ws.Range("L138").Value = 52800
$ws.Range("N138").Value = -63080

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 55144.5
$ws.Range("J57").Value = 55144.5
$ws.Range("L57").Value = 55144.5
$ws.Range("N57").Value = -56584.5
$ws.Range("H110").Value = 44600.75
$ws.Range("J110").Value = 44600.75
$ws.Range("L110").Value = 44600.75
$ws.Range("N110").Value = -52780.75
$ws.Range("H112").Value = 47469
$ws.Range("J112").Value = 47469
$ws.Range("L112").Value = 47469
$ws.Range("N112").Value = -50423
$ws.Range("H125").Value = 49772
$ws.Range("J125").Value = 49772
$ws.Range("L125").Value = 49772
$ws.Range("N125").Value = -59612
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656
$ws.Range("H133").Value = 48300
$ws.Range("J133").Value = 48300
$ws.Range("L133").Value = 48300
$ws.Range("N133").Value = -58420
$ws.Range("H136").Value = 55144.5
$ws.Range("J136").Value = 55144.5
$ws.Range("L136").Value = 55144.5
$ws.Range("N136").Value = -65344.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 65533
$ws.Range("J52").Value = 65533
$ws.Range("L52").Value = 65533
$ws.Range("N52").Value = -66121
$ws.Range("H110").Value = 40423
$ws.Range("J110").Value = 40423
$ws.Range("L110").Value = 40423
$ws.Range("N110").Value = -48603
$ws.Range("H112").Value = 47702
$ws.Range("J112").Value = 47702
$ws.Range("L112").Value = 47702
$ws.Range("N112").Value = -50656
$ws.Range("H137").Value = 35076.92
$ws.Range("J137").Value = 35076.92
$ws.Range("L137").Value = 35076.92
$ws.Range("N137").Value = -45276.92
$ws.Range("H139").Value = 59079.8
$ws.Range("J139").Value = 63349.75
$ws.Range("L139").Value = 63349.75
$ws.Range("N139").Value = -73629.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47011.332
$ws.Range("J110").Value = 47011.332
$ws.Range("L110").Value = 47011.332
$ws.Range("N110").Value = -55191.332
$ws.Range("H113").Value = 1251.5625
$ws.Range("J113").Value = 1260.4286
$ws.Range("L113").Value = 1260.4286
$ws.Range("N113").Value = -5600.4286
$ws.Range("H119").Value = 48437.5
$ws.Range("J119").Value = 48437.5
$ws.Range("L119").Value = 48437.5
$ws.Range("N119").Value = -58113.5
$ws.Range("H122").Value = 1257.5714
$ws.Range("I122").Value = 1257.5714
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3772.7142
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1322.7142
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H130").Value = 53986
$ws.Range("J130").Value = 53986
$ws.Range("L130").Value = 53986
$ws.Range("N130").Value = -64026
$ws.Range("H132").Value = 2806.6177
$ws.Range("I132").Value = 2031.65
$ws.Range("J132").Value = 3913.7144
$ws.Range("K132").Value = 6094.950000000001
$ws.Range("L132").Value = 11741.1432
$ws.Range("M132").Value = -3564.950000000001
$ws.Range("N132").Value = -16801.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 27188
$ws.Range("I2").Value = 1500
$ws.Range("K2").Value = 1500
$ws.Range("M2").Value = -1388
$ws.Range("H36").Value = 48460.75
$ws.Range("J36").Value = 48460.75
$ws.Range("L36").Value = 48460.75
$ws.Range("N36").Value = -49584.75
$ws.Range("H40").Value = 3641.5
$ws.Range("I40").Value = 2735.2856
$ws.Range("J40").Value = 9985
$ws.Range("K40").Value = 2735.2856
$ws.Range("L40").Value = 9985
$ws.Range("M40").Value = -2599.2856
$ws.Range("N40").Value = -10257
$ws.Range("H122").Value = 47669
$ws.Range("I122").Value = 73429.5
$ws.Range("J122").Value = 2588.125
$ws.Range("K122").Value = 220288.5
$ws.Range("L122").Value = 7764.375
$ws.Range("M122").Value = -217838.5
$ws.Range("N122").Value = -12664.375
$ws.Range("H124").Value = 47976
$ws.Range("J124").Value = 47976
$ws.Range("L124").Value = 47976
$ws.Range("N124").Value = -57796
$ws.Range("H125").Value = 49032.25
$ws.Range("J125").Value = 49032.25
$ws.Range("L125").Value = 49032.25
$ws.Range("N125").Value = -58872.25
$ws.Range("H127").Value = 48992
$ws.Range("J127").Value = 48992
$ws.Range("L127").Value = 48992
$ws.Range("N127").Value = -58912
$ws.Range("H128").Value = 32429
$ws.Range("J128").Value = 32429
$ws.Range("L128").Value = 32429
$ws.Range("N128").Value = -42389
$ws.Range("H130").Value = 48292
$ws.Range("J130").Value = 48292
$ws.Range("L130").Value = 48292
$ws.Range("N130").Value = -58332
$ws.Range("H137").Value = 40716.668
$ws.Range("J137").Value = 40716.668
$ws.Range("L137").Value = 40716.668
$ws.Range("N137").Value = -50916.668
$ws.Range("H139").Value = 68483
$ws.Range("I139").Value = 160000
$ws.Range("J139").Value = 50179.6
$ws.Range("K139").Value = 160000
$ws.Range("L139").Value = 50179.6
$ws.Range("M139").Value = -154860
$ws.Range("N139").Value = -60459.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 46671
$ws.Range("J46").Value = 46671
$ws.Range("L46").Value = 46671
$ws.Range("N46").Value = -47133
$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824
$ws.Range("H117").Value = 42571
$ws.Range("J117").Value = 42571
$ws.Range("L117").Value = 42571
$ws.Range("N117").Value = -51749
$ws.Range("H128").Value = 47984
$ws.Range("J128").Value = 47984
$ws.Range("L128").Value = 47984
$ws.Range("N128").Value = -57944
$ws.Range("H131").Value = 48982.5
$ws.Range("J131").Value = 48982.5
$ws.Range("L131").Value = 48982.5
$ws.Range("N131").Value = -59062.5
$ws.Range("H134").Value = 46671
$ws.Range("J134").Value = 46671
$ws.Range("L134").Value = 140013
$ws.Range("N134").Value = -145083
$ws.Range("H139").Value = 50919.8
$ws.Range("J139").Value = 50919.8
$ws.Range("L139").Value = 50919.8
$ws.Range("N139").Value = -61199.8
